# Atualização automática: 2025-08-12 09:01:14
# Updates two detection rows (16 and 17) with corrected image filename,
# bounding-box coordinates, and confidence score.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    # Force the cell to be treated as text so numeric-looking strings
    # (e.g. "641,530,687,575" or "0.76") are not coerced into numbers,
    # then restore the default "Normal" style so no stray number-format
    # style is left behind on the cell.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 16
Set-TextValue $ws.Range("D16") "image_20250807110238_ppp0.jpg"
Set-TextValue $ws.Range("I16") "641,530,687,575"
Set-TextValue $ws.Range("J16") "0.76"

# Row 17
Set-TextValue $ws.Range("D17") "image_20250807110238_ppp0.jpg"
Set-TextValue $ws.Range("I17") "793,481,831,527"
Set-TextValue $ws.Range("J17") "0.71"
